$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add notes: mark "Basics", "Control Structures" and "Functions" rows
# (JavaScript section) as completed by putting 1 in column E
$ws.Range("E20").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("E22").Value = 1

# Leave the view scrolled/selected where the user ended up editing
$ws.Activate()
$excel.Goto($ws.Range("A17"), $true)
$ws.Range("E23").Select()
